$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume/% change (E) columns with the latest scrape.
# Price cells are plain text in this sheet; for new values that look like
# ordinary numbers we briefly force a Text format so Excel keeps the literal
# string (preserving trailing zeros / decimal precision) instead of parsing
# it as a number, then restore the default "Normal" cell style so no stray
# formatting is left behind.

$ws.Cells.Item(2, 4).Value = "61.025.95"
$ws.Cells.Item(2, 5).Value = "  -0.43%  "
$ws.Cells.Item(3, 4).Value = "3.365.42"
$ws.Cells.Item(3, 5).Value = "  +2.23%  "
$ws.Cells.Item(4, 5).Value = "  +0.06%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "570.21"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  +1.01%  "
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "135.63"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  +8.19%  "
$ws.Cells.Item(7, 5).Value = "  +0.05%  "
$ws.Cells.Item(8, 4).Value = "3.363.61"
$ws.Cells.Item(8, 5).Value = "  +2.01%  "
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.475"
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).Value = "  +0.52%  "
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "7.61"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  +6.11%  "
$ws.Cells.Item(11, 5).Value = "  +4.86%  "
$ws.Cells.Item(12, 5).Value = "  +4.85%  "
$ws.Cells.Item(13, 4).Value = "3.941.36"
$ws.Cells.Item(13, 5).Value = "  +2.57%  "
$ws.Cells.Item(14, 5).Value = "  +2.59%  "
$ws.Cells.Item(15, 5).Value = "  +2.61%  "
$ws.Cells.Item(16, 4).Value = "3.358.44"
$ws.Cells.Item(16, 5).Value = "  +1.99%  "
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "25.10"
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 5).Value = "  +3.18%  "
$ws.Cells.Item(18, 4).Value = "61.101.97"
$ws.Cells.Item(18, 5).Value = "  -0.34%  "
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "13.98"
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  +7.52%  "
$ws.Cells.Item(20, 5).Value = "  +4.41%  "
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "9.36"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  +3.07%  "
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "374.84"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  +5.49%  "
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "0.572"
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  +4.80%  "
$ws.Cells.Item(24, 4).Value = "3.498.48"
$ws.Cells.Item(24, 5).Value = "  +2.52%  "
$ws.Cells.Item(25, 5).Value = "  +0.15%  "
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "70.74"
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  +1.29%  "
$ws.Cells.Item(27, 5).Value = "  +11.74%  "
$ws.Cells.Item(28, 5).Value = "  +14.62%  "
$ws.Cells.Item(29, 5).Value = "  +9.70%  "
$ws.Cells.Item(30, 5).Value = "  -0.58%  "
$ws.Cells.Item(31, 5).Value = "  +3.69%  "
$ws.Cells.Item(32, 5).Value = "  +5.70%  "
$ws.Cells.Item(33, 5).Value = "  +2.59%  "
$ws.Cells.Item(34, 5).Value = "  -0.03%  "
$ws.Cells.Item(35, 4).Value = "3.396.14"
$ws.Cells.Item(35, 5).Value = "  +2.61%  "
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "23.45"
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  +5.50%  "
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "5.54"
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  +2.80%  "
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "6.95"
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  +5.90%  "
$ws.Cells.Item(39, 5).Value = "  +5.08%  "
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "164.10"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  -1.11%  "
$ws.Cells.Item(41, 5).Value = "  +5.46%  "
$ws.Cells.Item(42, 5).Value = "  +0.25%  "
$ws.Cells.Item(43, 5).Value = "  +1.73%  "
$ws.Cells.Item(44, 5).Value = "  +5.97%  "
$ws.Cells.Item(45, 5).Value = "  +9.37%  "
$ws.Cells.Item(46, 5).Value = "  +1.77%  "
$ws.Cells.Item(47, 5).Value = "  +6.28%  "
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "23.04"
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  +3.62%  "
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "6.96"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  +6.54%  "
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "23.03"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  +14.02%  "
$ws.Cells.Item(51, 5).Value = "  +16.34%  "
